$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Indonesia): fix the ID, correct the Name typo, correct the Phone Code ---
$ws.Range("A2").Value = "21f66d32-6a6f-482f-88cb-2a3657251185"
$ws.Range("B2").Value = "Indonesia"
# "+62" looks numeric, so type it with a leading quote to force text, then drop
# the quote-prefix formatting so the cell keeps the default style.
$ws.Range("C2").Value = "'+62"
$ws.Range("C2").Style = "Normal"

# --- Row 4 (new): Malaysia / "Malasya" ---
$ws.Range("A4").Value = "7f7de98c-755e-46cc-ac65-cbf96b2204de"
$ws.Range("B4").Value = "Malasya"
$ws.Range("C4").Value = "'+63"
$ws.Range("C4").Style = "Normal"
# Icon Flag Path is blank, same as the other rows (stored as an empty text value).
$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = "Normal"

# --- Column widths: column B:C (width 12) splits into B=11, C=12 ---
# ColumnWidth (chars) -> stored xml width is ColumnWidth + 5/6, so subtract 5/6.
$ws.Columns.Item(2).ColumnWidth = 11 - 5/6
$ws.Columns.Item(3).ColumnWidth = 12 - 5/6
